$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6064498424530029
$ws.Range("B1").Value = 0.4509337544441223
$ws.Range("C1").Value = 0.4931572377681732
$ws.Range("D1").Value = 3.898004293441772
$ws.Range("E1").Value = 1.609934568405151
